# Edit script: reorders three pairs of match rows (6/7, 24/25, 30/31) so that
# their "home/away/odds/url" data (columns F:V) is swapped between the two
# rows in each pair, and appends two new match rows (42/43) at the end of
# the table, extending the used range to A1:V43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: swap the contents of a single column between two rows.
# Uses .Value2 for reading (reliable for both text and numeric cells in this
# engine) and .Value for writing.
# ---------------------------------------------------------------------------
function Swap-Cell($col, $r1, $r2) {
    $a = $ws.Range("$col$r1").Value2
    $b = $ws.Range("$col$r2").Value2
    $ws.Range("$col$r1").Value = $b
    $ws.Range("$col$r2").Value = $a
}

$dataCols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$swapPairs = @(
    @(6, 7),
    @(24, 25),
    @(30, 31)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($c in $dataCols) {
        Swap-Cell $c $r1 $r2
    }
}

# ---------------------------------------------------------------------------
# Append two new rows (42 and 43) after the current last row (41).
# Copy formatting from row 41 first so the new rows inherit the same styles
# (bold/centered index column, date-formatted match-date column) instead of
# creating brand-new style entries.
# ---------------------------------------------------------------------------
$ws.Range("A41:V41").Copy($ws.Range("A42:V42"))
$ws.Range("A41:V41").Copy($ws.Range("A43:V43"))

# Row 42
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = "italy"
$ws.Range("C42").Value = "serie-a"
$ws.Range("D42").Value = "2023-2024"
$ws.Range("E42").Value = 45191.77083333334
$ws.Range("F42").Value = "Salernitana"
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = "Frosinone"
$ws.Range("I42").Value = 1
$ws.Range("J42").Value = 2.41
$ws.Range("K42").Value = "11/09/2023 13:22"
$ws.Range("L42").Value = 2.45
$ws.Range("M42").Value = "22/09/2023 18:26"
$ws.Range("N42").Value = 3.16
$ws.Range("O42").Value = "11/09/2023 13:22"
$ws.Range("P42").Value = 3.36
$ws.Range("Q42").Value = "22/09/2023 18:20"
$ws.Range("R42").Value = 3.09
$ws.Range("S42").Value = "11/09/2023 13:22"
$ws.Range("T42").Value = 3.17
$ws.Range("U42").Value = "22/09/2023 18:28"
$ws.Range("V42").Value = "https://www.betexplorer.com/football/italy/serie-a/salernitana-frosinone/0KlD5LH4/"

# Row 43
$ws.Range("A43").Value = 42
$ws.Range("B43").Value = "italy"
$ws.Range("C43").Value = "serie-a"
$ws.Range("D43").Value = "2023-2024"
$ws.Range("E43").Value = 45191.86458333334
$ws.Range("F43").Value = "Lecce"
$ws.Range("G43").Value = 1
$ws.Range("H43").Value = "Genoa"
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2.1
$ws.Range("K43").Value = "11/09/2023 13:23"
$ws.Range("L43").Value = 2.65
$ws.Range("M43").Value = "22/09/2023 20:43"
$ws.Range("N43").Value = 3.17
$ws.Range("O43").Value = "11/09/2023 13:23"
$ws.Range("P43").Value = 3.19
$ws.Range("Q43").Value = "22/09/2023 19:58"
$ws.Range("R43").Value = 3.79
$ws.Range("S43").Value = "11/09/2023 13:23"
$ws.Range("T43").Value = 3.02
$ws.Range("U43").Value = "22/09/2023 20:43"
$ws.Range("V43").Value = "https://www.betexplorer.com/football/italy/serie-a/lecce-genoa/v5047smh/"
